$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed cryptos snapshot.
# Numeric-looking Price values need NumberFormat "@" forced beforehand so Excel
# keeps storing them as text (matching the source data), then the cell style is
# restored to Normal so no stray formatting is introduced.

$ws.Range("D2").Value = '34.738.83'
$ws.Range("E2").Value = '  -1.69%  '
$ws.Range("D3").Value = '1.809.70'
$ws.Range("E3").Value = '  -2.00%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.14%  '
$ws.Range("E6").Value = '  -1.45%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '39.30'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.91%  '
$ws.Range("E9").Value = '  +4.68%  '
$ws.Range("E10").Value = '  -1.60%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0993'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.46%  '
$ws.Range("D12").Value = '2.070.19'
$ws.Range("E12").Value = '  -2.07%  '
$ws.Range("D13").Value = '1.800.05'
$ws.Range("E13").Value = '  -2.47%  '
$ws.Range("E14").Value = '  -0.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.99'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.99%  '
$ws.Range("E16").Value = '  -1.92%  '
$ws.Range("D17").Value = '34.712.55'
$ws.Range("E17").Value = '  -1.76%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.48%  '
$ws.Range("D19").Value = '0.0₃0787'
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.64'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.20%  '
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("E24").Value = '  +2.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '172.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("E26").Value = '  -1.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.20'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.120'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +12.47%  '
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("E31").Value = '  +2.84%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0549'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E33").Value = '  -1.70%  '
$ws.Range("E34").Value = '  +17.00%  '
$ws.Range("E35").Value = '  -4.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.702'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.53%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '91.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.50%  '
$ws.Range("E38").Value = '  +5.86%  '
$ws.Range("D39").Value = '1.319.12'
$ws.Range("E39").Value = '  -2.09%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0192'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.12%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.965'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '14.21'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.44%  '
$ws.Range("E44").Value = '  -8.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.65'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.56%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.15%  '
$ws.Range("E47").Value = '  -1.58%  '
$ws.Range("D48").Value = '1.997.52'
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("E50").Value = '  +7.95%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '98.91'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.36%  '
